# CoinScales.xlsx - "Dynamic Token Models to be obtain via plist"
# Tokens moved to a circular shape for easier shape mapping:
# add a second, scaled-down (half-size) copy of the token-diameter table
# in rows 8-11 so the circular token shapes can be mapped to both the
# full-size and half-size values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: half of the row-3 (full-size) diameters.
$ws.Range("D8").Formula = "=D3/2"
$ws.Range("E8:F8").Formula = "=E3/2"

# Rows 9-11: half of the row 4-6 (Nickel/Dime/Penny) diameters, filled
# as one block so Excel groups them into a single shared formula - same
# shape the original D5:D6 "Dime" shared formula uses.
$ws.Range("D9:F11").Formula = "=D4/2"

# Keep the page set to portrait (matches the printed layout once the
# sheet grew past the original A3:F6 block).
$ws.PageSetup.Orientation = 1

# Leave the same cell selected as when the sheet was last saved.
$ws.Range("F9").Select() | Out-Null
